$d = $word.ActiveDocument

# Locate the introductory paragraph (begins "Ви сте учесници ...") that
# currently contains many differently-formatted runs about the
# "Персеус" constellation, and rewrite it as a single, unformatted run
# that references the "Близанци" (Gemini) constellation instead.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Ви сте учесници*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $newText = "Ви сте учесници глобалног посматрачког пројекта, који има за циљ да одреди колико је светлосно загађене у средини у којој живите. Посматрајући звезде унутар  сазвежђе Близанци и упоређујући их са приложеним звезданим картама, посматрачи широм света могу на практичном примеру да увиде колико је светлосно загађење у њиховој средини. Кроз учешће у овом пројекту, допринећете целовитијем сагледавању глобалног проблема."

    $start = $target.Range.Start
    $end = $target.Range.End - 1   # exclude the trailing paragraph mark
    $r = $d.Range($start, $end)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xml)
}
